# The commit adds one new daily price record for "Pepino ensalada" at
# Vega Central Mapocho de Santiago, inserted right before the existing
# row 344 (whose date is 2023-03-30 / serial 45015). Inserting the row
# pushes the former rows 344-409 down to 345-410, which is exactly the
# shift the diff shows (every row's data now equals the row above's old
# data, and a brand-new row 410 appears carrying the former row 409's
# values). The new row reuses the same market/category/quality/price
# data as the row it sits above, only the date changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 344; Excel shifts rows 344:409 down to 345:410
# and copies the row-above formatting (keeps the date style on column D).
$ws.Rows.Item(344).Insert()

$ws.Cells.Item(344, 1).Value = 9
$ws.Cells.Item(344, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(344, 3).Value = "Metropolitana"
$ws.Cells.Item(344, 4).Value = 45015
$ws.Cells.Item(344, 5).Value = 13
$ws.Cells.Item(344, 6).Value = 100112043
$ws.Cells.Item(344, 7).Value = "Pepino ensalada"
$ws.Cells.Item(344, 8).Value = "Sin especificar"
$ws.Cells.Item(344, 9).Value = "Primera"
$ws.Cells.Item(344, 10).Value = 70
$ws.Cells.Item(344, 11).Value = 7000
$ws.Cells.Item(344, 12).Value = 8000
$ws.Cells.Item(344, 13).Value = 7500
$ws.Cells.Item(344, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(344, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(344, 16).Value = 125
$ws.Cells.Item(344, 17).Value = 60
$ws.Cells.Item(344, 18).Value = "Hortaliza"
